$d = $word.ActiveDocument

# Locate the statement-of-truth sentence and replace it with the updated
# wording. The trailing "." is re-inserted as its own run (matching the
# target OOXML, which splits the sentence into two runs that share the
# same run formatting).
$rng = $d.Content
$found = $rng.Find.Execute(
    "The claimant believes that the facts in this claim form are true.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Text = "The Claimant believes that the facts stated in the brief details of claim are true"
    $rng.Collapse(0)
    $rng.InsertAfter(".")
    # Toggling formatting forces the inserted "." to stay in its own run
    # instead of being silently re-merged with the preceding run.
    $rng.Bold = 1
    $rng.Bold = 0
}
